$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every existing date_bin label down by one bin (a new earliest
# bin is being inserted), keeping the count column (A) untouched.
$ws.Range("B2").Value = "(2018-12-26, 2019-01-02]"
$ws.Range("B3").Value = "(2019-01-02, 2019-01-09]"
$ws.Range("B4").Value = "(2019-01-09, 2019-01-16]"
$ws.Range("B5").Value = "(2019-01-16, 2019-01-23]"
$ws.Range("B6").Value = "(2019-01-23, 2019-01-30]"
$ws.Range("B7").Value = "(2019-01-30, 2019-02-06]"
$ws.Range("B8").Value = "(2019-02-06, 2019-02-13]"
$ws.Range("B9").Value = "(2019-02-13, 2019-02-20]"
$ws.Range("B10").Value = "(2019-02-20, 2019-02-27]"
$ws.Range("B11").Value = "(2019-02-27, 2019-03-06]"
$ws.Range("B12").Value = "(2019-03-06, 2019-03-13]"
$ws.Range("B13").Value = "(2019-03-13, 2019-03-20]"
$ws.Range("B14").Value = "(2019-03-20, 2019-03-27]"
$ws.Range("B15").Value = "(2019-03-27, 2019-04-03]"
$ws.Range("B16").Value = "(2019-04-03, 2019-04-10]"
$ws.Range("B17").Value = "(2019-04-10, 2019-04-17]"
$ws.Range("B18").Value = "(2019-04-17, 2019-04-24]"
$ws.Range("B19").Value = "(2019-04-24, 2019-05-01]"
$ws.Range("B20").Value = "(2019-05-01, 2019-05-08]"
$ws.Range("B21").Value = "(2019-05-08, 2019-05-15]"
$ws.Range("B22").Value = "(2019-05-15, 2019-05-22]"
$ws.Range("B23").Value = "(2019-05-22, 2019-05-29]"
$ws.Range("B24").Value = "(2019-05-29, 2019-06-05]"
$ws.Range("B25").Value = "(2019-06-05, 2019-06-12]"
$ws.Range("B26").Value = "(2019-06-12, 2019-06-19]"
$ws.Range("B27").Value = "(2019-06-19, 2019-06-26]"
$ws.Range("B28").Value = "(2019-06-26, 2019-07-03]"
$ws.Range("B29").Value = "(2019-07-03, 2019-07-10]"
$ws.Range("B30").Value = "(2019-07-10, 2019-07-17]"
$ws.Range("B31").Value = "(2019-07-17, 2019-07-24]"
$ws.Range("B32").Value = "(2019-07-24, 2019-07-31]"
$ws.Range("B33").Value = "(2019-07-31, 2019-08-07]"
$ws.Range("B34").Value = "(2019-08-07, 2019-08-14]"
$ws.Range("B35").Value = "(2019-08-14, 2019-08-21]"
$ws.Range("B36").Value = "(2019-08-21, 2019-08-28]"
$ws.Range("B37").Value = "(2019-08-28, 2019-09-04]"
$ws.Range("B38").Value = "(2019-09-04, 2019-09-11]"
$ws.Range("B39").Value = "(2019-09-11, 2019-09-18]"
$ws.Range("B40").Value = "(2019-09-18, 2019-09-25]"
$ws.Range("B41").Value = "(2019-09-25, 2019-10-02]"
$ws.Range("B42").Value = "(2019-10-02, 2019-10-09]"
$ws.Range("B43").Value = "(2019-10-09, 2019-10-16]"
$ws.Range("B44").Value = "(2019-10-16, 2019-10-23]"
$ws.Range("B45").Value = "(2019-10-23, 2019-10-30]"
$ws.Range("B46").Value = "(2019-10-30, 2019-11-06]"
$ws.Range("B47").Value = "(2019-11-06, 2019-11-13]"
$ws.Range("B48").Value = "(2019-11-13, 2019-11-20]"
$ws.Range("B49").Value = "(2019-11-20, 2019-11-27]"
$ws.Range("B50").Value = "(2019-11-27, 2019-12-04]"
$ws.Range("B51").Value = "(2019-12-04, 2019-12-11]"

# Append the new row that re-uses the previous last bin label.
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "(2019-12-11, 2019-12-18]"
$ws.Range("A52").Style = $ws.Range("A51").Style
